$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("elifeChangePlanPage")

# The sheet currently has a duplicated "termsAndConditions" key in rows 21 & 22.
# Remove the duplicate row (row 22) and shift everything below it up by one,
# matching the updated "elife change plan" key list (now A1:B34 instead of A1:B35).
$ws.Rows.Item(22).Delete()

# Update the view state left behind by the edit (scrolled down, new active cell).
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("I34").Select()
